# M07 Frozen Token Embeddings + Decoder 1
#
# This re-run produced new per-epoch accuracy figures. Refresh column B
# (the accuracy column) for every row whose value actually moved, refresh
# the "<__main__.DisplayOutputs ...>" object-repr labels in column A whose
# embedded memory address changed between runs (A102:A118), and restore
# the sheet selection to the full used range (A2:B118), anchored at A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new accuracy value, for every B-column cell that changed.
$newAccuracy = @{
    2=0.921875; 4=0.875; 5=0.8125; 6=0.75; 9=0.640625; 10=0.640625; 11=0.640625; 12=0.65625;
    13=0.609375; 14=0.625; 15=0.640625; 16=0.5625; 17=0.515625; 18=0.609375; 19=0.453125; 20=0.515625;
    21=0.4375; 23=0.4375; 24=0.4375; 25=0.453125; 26=0.453125; 28=0.5; 29=0.5; 30=0.515625;
    31=0.515625; 32=0.515625; 33=0.515625; 34=0.515625; 35=0.515625; 36=0.515625; 37=0.515625; 38=0.515625;
    39=0.515625; 40=0.515625; 41=0.515625; 42=0.515625; 43=0.515625; 44=0.515625; 45=0.515625; 46=0.515625;
    47=0.515625; 48=0.515625; 49=0.515625; 50=0.515625; 51=0.515625; 52=0.515625; 53=0.515625; 54=0.515625;
    55=0.515625; 56=0.515625; 57=0.5; 58=0.5; 59=0.5; 60=0.5; 61=0.5; 62=0.5;
    63=0.5; 64=0.5; 65=0.5; 66=0.5; 67=0.5; 68=0.5; 69=0.5; 70=0.5;
    71=0.5; 72=0.5; 73=0.5; 74=0.5; 75=0.5; 76=0.5; 77=0.5; 78=0.5;
    79=0.5; 80=0.5; 81=0.5; 82=0.5; 83=0.5; 84=0.5; 85=0.5; 86=0.5;
    87=0.5; 88=0.5; 89=0.5; 90=0.5; 91=0.5; 92=0.5; 93=0.5; 94=0.5;
    95=0.5; 96=0.5; 97=0.5; 98=0.5; 99=0.5; 100=0.5; 101=0.5; 102=0.5;
    103=0.515625; 104=0.546875; 105=0.5; 106=0.484375; 107=0.5; 108=0.59375; 109=0.578125; 110=0.53125;
    111=0.5625; 113=0.5625; 114=0.5625; 115=0.484375; 116=0.5; 117=0.53125; 118=0.459016393442623
}

foreach ($entry in $newAccuracy.GetEnumerator()) {
    $ws.Cells.Item([int]$entry.Key, 2).Value = $entry.Value
}

# The DisplayOutputs repr in column A (rows 102-118) carries the Python
# object's memory address; it changed to 0x7fa841765430 on this re-run.
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = "<__main__.DisplayOutputs object at 0x7fa841765430>"
}

# Restore the selection to the sheet's used range, anchored at A2.
$ws.Range("A2:B118").Select()
